# Update average_county_temperature (col I), worst_ashp_cop (col N), and
# best_ashp_cop (col O) for rows whose underlying county temperature data
# changed after the NOAA refresh. Values below are the new canonical
# figures; rows where N/O are blank have no COP cells in the source row
# (non-electrifiable options do not carry those columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows  = @(10, 11, 20, 21, 24, 25, 26, 27, 34, 35, 40, 41, 48, 49, 50, 51, 52, 53, 56, 65, 66, 69, 70, 71, 72, 73, 74, 75, 76, 78, 79, 80, 81, 88, 89, 90, 91, 92, 93, 102, 103, 104, 105, 110, 111, 112, 113)
$ivals = @(13.62268518518517, 13.62268518518517, 12.93898809523811, 12.93898809523811, 13.46442495126706, 13.46442495126706, 15.74228395061728, 15.74228395061728, 19.79629629629628, 19.79629629629628, 20.68981481481483, 20.68981481481483, 21.28240740740739, 21.28240740740739, 19.65277777777778, 19.65277777777778, 21.28240740740739, 21.28240740740739, 19.65277777777778, 13.62268518518517, 13.62268518518517, 19.79629629629628, 19.79629629629628, 13.75752314814816, 13.75752314814816, 19.30324074074072, 19.30324074074072, 13.00385802469133, 13.00385802469133, 13.62268518518517, 13.62268518518517, 1.791666666666668, 1.791666666666668, 13.62268518518517, 13.62268518518517, 5.486111111111112, 5.486111111111112, 1.925925925925943, 1.925925925925943, 13.75752314814816, 13.75752314814816, 13.62268518518517, 13.62268518518517, 1.791666666666668, 1.791666666666668, 21.19907407407406, 21.19907407407406)
$nvals = @(0.9516859959956178, "", 0.9495669873270495, "", 0.9511946531936644, "", 0.9583158770332573, "", 0.9712571710566898, "", 0.9741566255742371, "", 0.9760891465058971, "", 0.970793063583815, "", 0.9760891465058971, "", "", 0.9516859959956178, "", 0.9712571710566898, "", 0.9521050214763401, "", 0.9696645907267841, "", 0.9497676359185355, "", 0.9516859959956178, "", 0.9163022129108289, "", 0.9516859959956178, "", 0.9270655773901523, "", 0.9166889846297158, "", 0.9521050214763401, "", 0.9516859959956178, "", 0.9163022129108289, "", 0.9758169225763468, "")
$ovals = @(0.9753614246104579, "", 0.9731165936130245, "", 0.9748408768576692, "", 0.9823871532785465, "", 0.9961106502456767, "", 0.9991870829399434, "", 1.001237913506406, "", 0.995618273140397, "", 1.001237913506406, "", "", 0.9753614246104579, "", 0.9961106502456767, "", 0.9758053708974481, "", 0.9944211305850406, "", 0.9733291418446532, "", 0.9753614246104579, "", 0.9379207786940652, "", 0.9753614246104579, "", 0.9492998859749143, "", 0.9383295263284442, "", 0.9758053708974481, "", 0.9753614246104579, "", 0.9379207786940652, "", 1.000949006909155, "")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 9).Value = $ivals[$i]   # column I

    $n = $nvals[$i]
    if ($n -ne "") {
        $ws.Cells.Item($r, 14).Value = $n      # column N
    }

    $o = $ovals[$i]
    if ($o -ne "") {
        $ws.Cells.Item($r, 15).Value = $o      # column O
    }
}
